$d = $word.ActiveDocument
$p5 = $d.Paragraphs(5)
$r = $p5.Range
$r.Collapse(1)
$r.InsertParagraphBefore()
$newPara = $d.Paragraphs(5)
$newPara.Range.Text = "Change paragraphs to boxes"
$newPara.Style = "Section3"

# Insert second new paragraph after newPara, before Portfolio/Projects (now at 6)
$p6 = $d.Paragraphs(6)
$r6 = $p6.Range
$r6.Collapse(1)
$r6.InsertParagraphBefore()
$newPara2 = $d.Paragraphs(6)
$newPara2.Style = "Section3"
Write-Host "newPara2 text: [" $newPara2.Range.Text "]"
$d.Bookmarks.Add("_GoBack", $newPara2.Range)
Write-Host "Done"
